$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every Price ("D") cell as literal text (inline
# strings), even when the text happens to look like a plain number (e.g.
# "1.004"). Mark the whole Price column as Text first, so the updated
# values keep the same text representation as before instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.083.19"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "2.104.72"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  -0.64%  "
$ws.Range("D5").Value = "347.86"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").Value = "0.5175"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("D8").Value = "0.4443"
$ws.Range("E8").Value = "  -2.57%  "
$ws.Range("D9").Value = "52.41"
$ws.Range("E9").Value = "  -3.70%  "
$ws.Range("D10").Value = "0.08965"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").Value = "1.175"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  +4.13%  "
$ws.Range("D13").Value = "2.108.44"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "8.267"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "6.733"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "99.52"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "0.00001151"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "20.79"
$ws.Range("E19").Value = "  +7.11%  "
$ws.Range("D20").Value = "0.06688"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "6.248"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "30.164.82"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").Value = "12.74"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").Value = "2.340"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").Value = "2.356.39"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "21.96"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("D28").Value = "2.539"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "162.26"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "133.50"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "1.175"
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").Value = "1.636"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").Value = "6.242"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").Value = "3.960"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "5.955"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").Value = "10.26"
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("D38").Value = "0.02579"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "0.06813"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "0.2297"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").Value = "12.61"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "0.6817"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").Value = "1.283"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("D44").Value = "14.31"
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("D45").Value = "0.6388"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").Value = "2.296"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "3.646"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").Value = "1.220"
$ws.Range("E49").Value = "  -2.67%  "
$ws.Range("D50").Value = "82.46"
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").Value = "0.07233"
$ws.Range("E51").Value = "  +0.48%  "

Write-Host "Updated cryptos list"
